$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like columns (Price, Volume) keep their literal string
# representation instead of being auto-converted to numbers by Excel,
# since values like "0.620", "1.00", "2.97" would otherwise lose
# trailing zeros / significant digits when coerced to a numeric type.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.312.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.153.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.86"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.64"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.42%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0852"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.02"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.473.36"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.23"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.816"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.53"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.149.99"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "39.408.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.08"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.14"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.03"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.75"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.03"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.139"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.63"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.42%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +9.79%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.82"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.13"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +11.31%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.16"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.82"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.533.87"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.28%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.14%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0920"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.17"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.357.67"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.97"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.01%  "
